$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Organic N extraction id (B4: 102 -> 104)
$ws.Range("B4").Value = 104

# Populate the previously-blank lab values (B18:B29, B32) pulled from the
# re-extracted soil analysis table
$ws.Range("B18").Value = 17
$ws.Range("B19").Value = 8
$ws.Range("B20").Value = 9
$ws.Range("B21").Value = 2
$ws.Range("B22").Value = 8
$ws.Range("B23").Value = 45
$ws.Range("B24").Value = 12
$ws.Range("B25").Value = 15
$ws.Range("B26").Value = 42
$ws.Range("B27").Value = 2
$ws.Range("B28").Value = 4
$ws.Range("B29").Value = 1
$ws.Range("B32").Value = 1

# Move/save the active selection on B4 (also clears the stale scrolled
# topLeftCell from the view)
$ws.Range("B4").Select()
